$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the AF column ("size"/"sz") with "x" for every vehicle data row (3-20).
$ws.Range("AF3").Value = "x"
$ws.Range("AF4").Value = "x"
$ws.Range("AF5").Value = "x"
$ws.Range("AF6").Value = "x"
$ws.Range("AF7").Value = "x"
$ws.Range("AF8").Value = "x"
$ws.Range("AF9").Value = "x"
$ws.Range("AF10").Value = "x"
$ws.Range("AF11").Value = "x"
$ws.Range("AF12").Value = "x"
$ws.Range("AF13").Value = "x"
$ws.Range("AF14").Value = "x"
$ws.Range("AF15").Value = "x"
$ws.Range("AF16").Value = "x"
$ws.Range("AF17").Value = "x"
$ws.Range("AF18").Value = "x"
$ws.Range("AF19").Value = "x"
$ws.Range("AF20").Value = "x"

# Fill in previously-blank "inventory drop code" (Y) cells for the treb/hwacha rows.
$ws.Range("Y15").Value = "x"
$ws.Range("Y16").Value = "x"
$ws.Range("Y17").Value = "x"
$ws.Range("Y18").Value = "x"
$ws.Range("Y19").Value = "x"

# Row 20 (chest cart) also gets an "AE" inventory/BB entry.
$ws.Range("AE20").Value = "x"

# Move the active selection to reflect where editing left off.
$ws.Range("Y20").Select()
